# Refresh the COVID-19 "Pais" dashboard with the latest figures.
# The source feed re-ranks countries by total cases (column B, desc),
# so besides updating several country counters this also reshuffles
# a number of rows into their new rank position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" banner (row 1)
$ws.Cells.Item(1, 1).Value2 = 'Datos actualizados a 27 de Marzo de 2020 a las 18:14'

# Full league table, rows 4-205: Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = New-Object "object[,]" 202,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 93410
$data[0,2] = 7975
$data[0,3] = 2424
$data[0,4] = 89601
$data[0,5] = 2432
$data[0,6] = 90
$data[0,7] = 1385
$data[1,0] = 'Italia'
$data[1,1] = 86498
$data[1,2] = 5909
$data[1,3] = 10950
$data[1,4] = 66414
$data[1,5] = 3732
$data[1,6] = 919
$data[1,7] = 9134
$data[2,0] = 'China'
$data[2,1] = 81340
$data[2,2] = 55
$data[2,3] = 74588
$data[2,4] = 3460
$data[2,5] = 1034
$data[2,6] = 5
$data[2,7] = 3292
$data[3,0] = 'España'
$data[3,1] = 64059
$data[3,2] = 6273
$data[3,3] = 9357
$data[3,4] = 49768
$data[3,5] = 4165
$data[3,6] = 569
$data[3,7] = 4934
$data[4,0] = 'Alemania'
$data[4,1] = 49344
$data[4,2] = 5406
$data[4,3] = 5673
$data[4,4] = 43367
$data[4,5] = 23
$data[4,6] = 37
$data[4,7] = 304
$data[5,0] = 'Iran'
$data[5,1] = 32332
$data[5,2] = 2926
$data[5,3] = 11133
$data[5,4] = 18821
$data[5,5] = 2893
$data[5,6] = 144
$data[5,7] = 2378
$data[6,0] = 'Francia'
$data[6,1] = 29155
$data[6,2] = 0
$data[6,3] = 4948
$data[6,4] = 22511
$data[6,5] = 3375
$data[6,6] = 0
$data[6,7] = 1696
$data[7,0] = 'Reino Unido'
$data[7,1] = 14543
$data[7,2] = 2885
$data[7,3] = 135
$data[7,4] = 13649
$data[7,5] = 163
$data[7,6] = 181
$data[7,7] = 759
$data[8,0] = 'Suiza'
$data[8,1] = 12311
$data[8,2] = 500
$data[8,3] = 897
$data[8,4] = 11207
$data[8,5] = 203
$data[8,6] = 15
$data[8,7] = 207
$data[9,0] = 'Corea del Sur'
$data[9,1] = 9332
$data[9,2] = 91
$data[9,3] = 4528
$data[9,4] = 4665
$data[9,5] = 59
$data[9,6] = 8
$data[9,7] = 139
$data[10,0] = 'Paises Bajos'
$data[10,1] = 8603
$data[10,2] = 1172
$data[10,3] = 3
$data[10,4] = 8054
$data[10,5] = 761
$data[10,6] = 112
$data[10,7] = 546
$data[11,0] = 'Austria'
$data[11,1] = 7557
$data[11,2] = 648
$data[11,3] = 225
$data[11,4] = 7274
$data[11,5] = 128
$data[11,6] = 9
$data[11,7] = 58
$data[12,0] = 'Belgica'
$data[12,1] = 7284
$data[12,2] = 1049
$data[12,3] = 858
$data[12,4] = 6137
$data[12,5] = 690
$data[12,6] = 69
$data[12,7] = 289
$data[13,0] = 'Turquia'
$data[13,1] = 5698
$data[13,2] = 2069
$data[13,3] = 42
$data[13,4] = 5564
$data[13,5] = 241
$data[13,6] = 17
$data[13,7] = 92
$data[14,0] = 'Portugal'
$data[14,1] = 4268
$data[14,2] = 724
$data[14,3] = 43
$data[14,4] = 4149
$data[14,5] = 71
$data[14,6] = 16
$data[14,7] = 76
$data[15,0] = 'Canada'
$data[15,1] = 4043
$data[15,2] = 0
$data[15,3] = 228
$data[15,4] = 3776
$data[15,5] = 120
$data[15,6] = 0
$data[15,7] = 39
$data[16,0] = 'Noruega'
$data[16,1] = 3696
$data[16,2] = 324
$data[16,3] = 6
$data[16,4] = 3671
$data[16,5] = 70
$data[16,6] = 5
$data[16,7] = 19
$data[17,0] = 'Australia'
$data[17,1] = 3180
$data[17,2] = 130
$data[17,3] = 170
$data[17,4] = 2997
$data[17,5] = 23
$data[17,6] = 0
$data[17,7] = 13
$data[18,0] = 'Suecia'
$data[18,1] = 3046
$data[18,2] = 206
$data[18,3] = 16
$data[18,4] = 2938
$data[18,5] = 214
$data[18,6] = 15
$data[18,7] = 92
$data[19,0] = 'Israel'
$data[19,1] = 3035
$data[19,2] = 342
$data[19,3] = 79
$data[19,4] = 2944
$data[19,5] = 49
$data[19,6] = 4
$data[19,7] = 12
$data[20,0] = 'Brasil'
$data[20,1] = 3027
$data[20,2] = 42
$data[20,3] = 6
$data[20,4] = 2944
$data[20,5] = 296
$data[20,6] = 0
$data[20,7] = 77
$data[21,0] = 'Malasia'
$data[21,1] = 2161
$data[21,2] = 130
$data[21,3] = 259
$data[21,4] = 1876
$data[21,5] = 54
$data[21,6] = 3
$data[21,7] = 26
$data[22,0] = 'Chequia'
$data[22,1] = 2062
$data[22,2] = 137
$data[22,3] = 11
$data[22,4] = 2042
$data[22,5] = 34
$data[22,6] = 0
$data[22,7] = 9
$data[23,0] = 'Dinamarca'
$data[23,1] = 2010
$data[23,2] = 133
$data[23,3] = 1
$data[23,4] = 1957
$data[23,5] = 109
$data[23,6] = 11
$data[23,7] = 52
$data[24,0] = 'Irlanda'
$data[24,1] = 1819
$data[24,2] = 0
$data[24,3] = 5
$data[24,4] = 1795
$data[24,5] = 47
$data[24,6] = 0
$data[24,7] = 19
$data[25,0] = 'Chile'
$data[25,1] = 1610
$data[25,2] = 304
$data[25,3] = 43
$data[25,4] = 1562
$data[25,5] = 7
$data[25,6] = 1
$data[25,7] = 5
$data[26,0] = 'Luxemburgo'
$data[26,1] = 1605
$data[26,2] = 152
$data[26,3] = 40
$data[26,4] = 1550
$data[26,5] = 25
$data[26,6] = 6
$data[26,7] = 15
$data[27,0] = 'Ecuador'
$data[27,1] = 1595
$data[27,2] = 192
$data[27,3] = 3
$data[27,4] = 1556
$data[27,5] = 58
$data[27,6] = 2
$data[27,7] = 36
$data[28,0] = 'Japon'
$data[28,1] = 1387
$data[28,2] = 0
$data[28,3] = 359
$data[28,4] = 981
$data[28,5] = 57
$data[28,6] = 0
$data[28,7] = 47
$data[29,0] = 'Polonia'
$data[29,1] = 1340
$data[29,2] = 119
$data[29,3] = 7
$data[29,4] = 1317
$data[29,5] = 3
$data[29,6] = 0
$data[29,7] = 16
$data[30,0] = 'Pakistan'
$data[30,1] = 1296
$data[30,2] = 95
$data[30,3] = 23
$data[30,4] = 1264
$data[30,5] = 7
$data[30,6] = 0
$data[30,7] = 9
$data[31,0] = 'Rumania'
$data[31,1] = 1292
$data[31,2] = 263
$data[31,3] = 115
$data[31,4] = 1152
$data[31,5] = 32
$data[31,6] = 2
$data[31,7] = 25
$data[32,0] = 'Sudafrica'
$data[32,1] = 1170
$data[32,2] = 243
$data[32,3] = 12
$data[32,4] = 1156
$data[32,5] = 7
$data[32,6] = 2
$data[32,7] = 2
$data[33,0] = 'Tailandia'
$data[33,1] = 1136
$data[33,2] = 91
$data[33,3] = 97
$data[33,4] = 1034
$data[33,5] = 11
$data[33,6] = 1
$data[33,7] = 5
$data[34,0] = 'Arabia Saudita'
$data[34,1] = 1104
$data[34,2] = 92
$data[34,3] = 35
$data[34,4] = 1066
$data[34,5] = 6
$data[34,6] = 0
$data[34,7] = 3
$data[35,0] = 'Indonesia'
$data[35,1] = 1046
$data[35,2] = 153
$data[35,3] = 46
$data[35,4] = 913
$data[35,5] = 0
$data[35,6] = 9
$data[35,7] = 87
$data[36,0] = 'Finlandia'
$data[36,1] = 1041
$data[36,2] = 83
$data[36,3] = 10
$data[36,4] = 1024
$data[36,5] = 32
$data[36,6] = 2
$data[36,7] = 7
$data[37,0] = 'Rusia'
$data[37,1] = 1036
$data[37,2] = 196
$data[37,3] = 45
$data[37,4] = 988
$data[37,5] = 8
$data[37,6] = 0
$data[37,7] = 3
$data[38,0] = 'Grecia'
$data[38,1] = 966
$data[38,2] = 74
$data[38,3] = 52
$data[38,4] = 886
$data[38,5] = 66
$data[38,6] = 1
$data[38,7] = 28
$data[39,0] = 'Islandia'
$data[39,1] = 890
$data[39,2] = 88
$data[39,3] = 97
$data[39,4] = 791
$data[39,5] = 18
$data[39,6] = 0
$data[39,7] = 2
$data[40,0] = 'India'
$data[40,1] = 878
$data[40,2] = 151
$data[40,3] = 73
$data[40,4] = 785
$data[40,5] = 0
$data[40,6] = 0
$data[40,7] = 20
$data[41,0] = 'Filipinas'
$data[41,1] = 803
$data[41,2] = 96
$data[41,3] = 31
$data[41,4] = 718
$data[41,5] = 1
$data[41,6] = 9
$data[41,7] = 54
$data[42,0] = 'Singapur'
$data[42,1] = 732
$data[42,2] = 49
$data[42,3] = 183
$data[42,4] = 547
$data[42,5] = 17
$data[42,6] = 0
$data[42,7] = 2
$data[43,0] = 'Crucero'
$data[43,1] = 712
$data[43,2] = 0
$data[43,3] = 597
$data[43,4] = 105
$data[43,5] = 15
$data[43,6] = 0
$data[43,7] = 10
$data[44,0] = 'Panama'
$data[44,1] = 674
$data[44,2] = 0
$data[44,3] = 2
$data[44,4] = 663
$data[44,5] = 20
$data[44,6] = 0
$data[44,7] = 9
$data[45,0] = 'Eslovenia'
$data[45,1] = 632
$data[45,2] = 70
$data[45,3] = 10
$data[45,4] = 613
$data[45,5] = 14
$data[45,6] = 3
$data[45,7] = 9
$data[46,0] = 'Argentina'
$data[46,1] = 589
$data[46,2] = 0
$data[46,3] = 72
$data[46,4] = 504
$data[46,5] = 0
$data[46,6] = 1
$data[46,7] = 13
$data[47,0] = 'Croacia'
$data[47,1] = 586
$data[47,2] = 91
$data[47,3] = 37
$data[47,4] = 546
$data[47,5] = 14
$data[47,6] = 0
$data[47,7] = 3
$data[48,0] = 'Mexico'
$data[48,1] = 585
$data[48,2] = 110
$data[48,3] = 4
$data[48,4] = 573
$data[48,5] = 1
$data[48,6] = 2
$data[48,7] = 8
$data[49,0] = 'Republica Dominicana'
$data[49,1] = 581
$data[49,2] = 93
$data[49,3] = 3
$data[49,4] = 558
$data[49,5] = 0
$data[49,6] = 10
$data[49,7] = 20
$data[50,0] = 'Peru'
$data[50,1] = 580
$data[50,2] = 0
$data[50,3] = 14
$data[50,4] = 557
$data[50,5] = 14
$data[50,6] = 0
$data[50,7] = 9
$data[51,0] = 'Estonia'
$data[51,1] = 575
$data[51,2] = 37
$data[51,3] = 11
$data[51,4] = 563
$data[51,5] = 7
$data[51,6] = 0
$data[51,7] = 1
$data[52,0] = 'Catar'
$data[52,1] = 562
$data[52,2] = 13
$data[52,3] = 43
$data[52,4] = 519
$data[52,5] = 6
$data[52,6] = 0
$data[52,7] = 0
$data[53,0] = 'Serbia'
$data[53,1] = 528
$data[53,2] = 71
$data[53,3] = 15
$data[53,4] = 505
$data[53,5] = 25
$data[53,6] = 1
$data[53,7] = 8
$data[54,0] = 'Hong Kong'
$data[54,1] = 518
$data[54,2] = 64
$data[54,3] = 111
$data[54,4] = 403
$data[54,5] = 5
$data[54,6] = 0
$data[54,7] = 4
$data[55,0] = 'Egipto'
$data[55,1] = 495
$data[55,2] = 0
$data[55,3] = 102
$data[55,4] = 369
$data[55,5] = 0
$data[55,6] = 0
$data[55,7] = 24
$data[56,0] = 'Colombia'
$data[56,1] = 491
$data[56,2] = 0
$data[56,3] = 8
$data[56,4] = 477
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 6
$data[57,0] = 'Barein'
$data[57,1] = 466
$data[57,2] = 8
$data[57,3] = 227
$data[57,4] = 235
$data[57,5] = 1
$data[57,6] = 0
$data[57,7] = 4
$data[58,0] = 'Irak'
$data[58,1] = 458
$data[58,2] = 76
$data[58,3] = 122
$data[58,4] = 296
$data[58,5] = 0
$data[58,6] = 4
$data[58,7] = 40
$data[59,0] = 'Argelia'
$data[59,1] = 409
$data[59,2] = 42
$data[59,3] = 29
$data[59,4] = 354
$data[59,5] = 0
$data[59,6] = 1
$data[59,7] = 26
$data[60,0] = 'Libano'
$data[60,1] = 391
$data[60,2] = 23
$data[60,3] = 23
$data[60,4] = 361
$data[60,5] = 3
$data[60,6] = 1
$data[60,7] = 7
$data[61,0] = 'Nueva Zelanda'
$data[61,1] = 368
$data[61,2] = 0
$data[61,3] = 37
$data[61,4] = 331
$data[61,5] = 1
$data[61,6] = 0
$data[61,7] = 0
$data[62,0] = 'Lituania'
$data[62,1] = 345
$data[62,2] = 46
$data[62,3] = 1
$data[62,4] = 339
$data[62,5] = 2
$data[62,6] = 1
$data[62,7] = 5
$data[63,0] = 'Marruecos'
$data[63,1] = 333
$data[63,2] = 58
$data[63,3] = 8
$data[63,4] = 314
$data[63,5] = 1
$data[63,6] = 0
$data[63,7] = 11
$data[64,0] = 'Emiratos Arabes Unidos'
$data[64,1] = 333
$data[64,2] = 0
$data[64,3] = 52
$data[64,4] = 279
$data[64,5] = 2
$data[64,6] = 0
$data[64,7] = 2
$data[65,0] = 'Armenia'
$data[65,1] = 329
$data[65,2] = 39
$data[65,3] = 28
$data[65,4] = 300
$data[65,5] = 6
$data[65,6] = 0
$data[65,7] = 1
$data[66,0] = 'Hungria'
$data[66,1] = 300
$data[66,2] = 39
$data[66,3] = 34
$data[66,4] = 256
$data[66,5] = 6
$data[66,6] = 0
$data[66,7] = 10
$data[67,0] = 'Bulgaria'
$data[67,1] = 293
$data[67,2] = 29
$data[67,3] = 9
$data[67,4] = 281
$data[67,5] = 8
$data[67,6] = 0
$data[67,7] = 3
$data[68,0] = 'Letonia'
$data[68,1] = 280
$data[68,2] = 36
$data[68,3] = 1
$data[68,4] = 279
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 0
$data[69,0] = 'Eslovaquia'
$data[69,1] = 269
$data[69,2] = 43
$data[69,3] = 2
$data[69,4] = 267
$data[69,5] = 1
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = 'Principado de Andorra'
$data[70,1] = 267
$data[70,2] = 43
$data[70,3] = 1
$data[70,4] = 263
$data[70,5] = 11
$data[70,6] = 0
$data[70,7] = 3
$data[71,0] = 'Taiwan'
$data[71,1] = 267
$data[71,2] = 15
$data[71,3] = 30
$data[71,4] = 235
$data[71,5] = 0
$data[71,6] = 0
$data[71,7] = 2
$data[72,0] = 'Uruguay'
$data[72,1] = 238
$data[72,2] = 0
$data[72,3] = 0
$data[72,4] = 238
$data[72,5] = 3
$data[72,6] = 0
$data[72,7] = 0
$data[73,0] = 'Costa Rica'
$data[73,1] = 231
$data[73,2] = 0
$data[73,3] = 2
$data[73,4] = 227
$data[73,5] = 5
$data[73,6] = 0
$data[73,7] = 2
$data[74,0] = 'Bosnia y Herzegovina'
$data[74,1] = 231
$data[74,2] = 40
$data[74,3] = 5
$data[74,4] = 222
$data[74,5] = 1
$data[74,6] = 1
$data[74,7] = 4
$data[75,0] = 'Tunez'
$data[75,1] = 227
$data[75,2] = 30
$data[75,3] = 2
$data[75,4] = 219
$data[75,5] = 10
$data[75,6] = 1
$data[75,7] = 6
$data[76,0] = 'Ucrania'
$data[76,1] = 226
$data[76,2] = 30
$data[76,3] = 5
$data[76,4] = 216
$data[76,5] = 0
$data[76,6] = 0
$data[76,7] = 5
$data[77,0] = 'Kuwait'
$data[77,1] = 225
$data[77,2] = 17
$data[77,3] = 57
$data[77,4] = 168
$data[77,5] = 11
$data[77,6] = 0
$data[77,7] = 0
$data[78,0] = 'San Marino'
$data[78,1] = 223
$data[78,2] = 15
$data[78,3] = 4
$data[78,4] = 198
$data[78,5] = 12
$data[78,6] = 0
$data[78,7] = 21
$data[79,0] = 'Republica de Macedonia'
$data[79,1] = 219
$data[79,2] = 18
$data[79,3] = 3
$data[79,4] = 213
$data[79,5] = 1
$data[79,6] = 0
$data[79,7] = 3
$data[80,0] = 'Jordania'
$data[80,1] = 212
$data[80,2] = 0
$data[80,3] = 2
$data[80,4] = 210
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 0
$data[81,0] = 'Moldavia'
$data[81,1] = 199
$data[81,2] = 22
$data[81,3] = 2
$data[81,4] = 195
$data[81,5] = 33
$data[81,6] = 1
$data[81,7] = 2
$data[82,0] = 'Albania'
$data[82,1] = 186
$data[82,2] = 12
$data[82,3] = 31
$data[82,4] = 147
$data[82,5] = 3
$data[82,6] = 2
$data[82,7] = 8
$data[83,0] = 'Burkina Faso'
$data[83,1] = 180
$data[83,2] = 28
$data[83,3] = 12
$data[83,4] = 159
$data[83,5] = 0
$data[83,6] = 2
$data[83,7] = 9
$data[84,0] = 'Azerbaiyan'
$data[84,1] = 165
$data[84,2] = 43
$data[84,3] = 15
$data[84,4] = 147
$data[84,5] = 6
$data[84,6] = 0
$data[84,7] = 3
$data[85,0] = 'Vietnam'
$data[85,1] = 163
$data[85,2] = 10
$data[85,3] = 20
$data[85,4] = 143
$data[85,5] = 3
$data[85,6] = 0
$data[85,7] = 0
$data[86,0] = 'Republica de Chipre'
$data[86,1] = 146
$data[86,2] = 0
$data[86,3] = 4
$data[86,4] = 137
$data[86,5] = 3
$data[86,6] = 2
$data[86,7] = 5
$data[87,0] = 'Reunion'
$data[87,1] = 145
$data[87,2] = 10
$data[87,3] = 1
$data[87,4] = 144
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 0
$data[88,0] = 'Islas Feroe'
$data[88,1] = 144
$data[88,2] = 4
$data[88,3] = 54
$data[88,4] = 90
$data[88,5] = 2
$data[88,6] = 0
$data[88,7] = 0
$data[89,0] = 'Malta'
$data[89,1] = 139
$data[89,2] = 5
$data[89,3] = 2
$data[89,4] = 137
$data[89,5] = 1
$data[89,6] = 0
$data[89,7] = 0
$data[90,0] = 'Kazajistan'
$data[90,1] = 137
$data[90,2] = 24
$data[90,3] = 3
$data[90,4] = 133
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 1
$data[91,0] = 'Ghana'
$data[91,1] = 136
$data[91,2] = 4
$data[91,3] = 1
$data[91,4] = 131
$data[91,5] = 1
$data[91,6] = 0
$data[91,7] = 4
$data[92,0] = 'Oman'
$data[92,1] = 131
$data[92,2] = 22
$data[92,3] = 23
$data[92,4] = 108
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 0
$data[93,0] = 'Senegal'
$data[93,1] = 119
$data[93,2] = 14
$data[93,3] = 11
$data[93,4] = 108
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 0
$data[94,0] = 'Brunei'
$data[94,1] = 115
$data[94,2] = 1
$data[94,3] = 11
$data[94,4] = 104
$data[94,5] = 1
$data[94,6] = 0
$data[94,7] = 0
$data[95,0] = 'Venezuela'
$data[95,1] = 107
$data[95,2] = 0
$data[95,3] = 31
$data[95,4] = 75
$data[95,5] = 2
$data[95,6] = 0
$data[95,7] = 1
$data[96,0] = 'Sri Lanka'
$data[96,1] = 106
$data[96,2] = 0
$data[96,3] = 7
$data[96,4] = 99
$data[96,5] = 5
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = 'Camboya'
$data[97,1] = 99
$data[97,2] = 1
$data[97,3] = 11
$data[97,4] = 88
$data[97,5] = 1
$data[97,6] = 0
$data[97,7] = 0
$data[98,0] = 'Costa de Marfil'
$data[98,1] = 96
$data[98,2] = 0
$data[98,3] = 3
$data[98,4] = 93
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 0
$data[99,0] = 'Mauricio'
$data[99,1] = 94
$data[99,2] = 13
$data[99,3] = 0
$data[99,4] = 92
$data[99,5] = 1
$data[99,6] = 0
$data[99,7] = 2
$data[100,0] = 'Afganistan'
$data[100,1] = 94
$data[100,2] = 0
$data[100,3] = 2
$data[100,4] = 88
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 4
$data[101,0] = 'Bielorrusia'
$data[101,1] = 94
$data[101,2] = 8
$data[101,3] = 32
$data[101,4] = 62
$data[101,5] = 2
$data[101,6] = 0
$data[101,7] = 0
$data[102,0] = 'Estado de Palestina'
$data[102,1] = 91
$data[102,2] = 5
$data[102,3] = 17
$data[102,4] = 73
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 1
$data[103,0] = 'Camerun'
$data[103,1] = 88
$data[103,2] = 13
$data[103,3] = 2
$data[103,4] = 84
$data[103,5] = 0
$data[103,6] = 1
$data[103,7] = 2
$data[104,0] = 'Uzbekistan'
$data[104,1] = 88
$data[104,2] = 13
$data[104,3] = 5
$data[104,4] = 82
$data[104,5] = 8
$data[104,6] = 1
$data[104,7] = 1
$data[105,0] = 'Martinica'
$data[105,1] = 81
$data[105,2] = 0
$data[105,3] = 0
$data[105,4] = 80
$data[105,5] = 12
$data[105,6] = 0
$data[105,7] = 1
$data[106,0] = 'Georgia'
$data[106,1] = 81
$data[106,2] = 2
$data[106,3] = 13
$data[106,4] = 68
$data[106,5] = 1
$data[106,6] = 0
$data[106,7] = 0
$data[107,0] = 'Cuba'
$data[107,1] = 80
$data[107,2] = 13
$data[107,3] = 4
$data[107,4] = 74
$data[107,5] = 2
$data[107,6] = 0
$data[107,7] = 2
$data[108,0] = 'Montenegro'
$data[108,1] = 75
$data[108,2] = 6
$data[108,3] = 0
$data[108,4] = 74
$data[108,5] = 1
$data[108,6] = 0
$data[108,7] = 1
$data[109,0] = 'Guadalupe'
$data[109,1] = 73
$data[109,2] = 0
$data[109,3] = 0
$data[109,4] = 72
$data[109,5] = 4
$data[109,6] = 0
$data[109,7] = 1
$data[110,0] = 'Honduras'
$data[110,1] = 68
$data[110,2] = 1
$data[110,3] = 0
$data[110,4] = 67
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 1
$data[111,0] = 'Trinidad yTobago'
$data[111,1] = 66
$data[111,2] = 1
$data[111,3] = 1
$data[111,4] = 63
$data[111,5] = 0
$data[111,6] = 1
$data[111,7] = 2
$data[112,0] = 'Nigeria'
$data[112,1] = 65
$data[112,2] = 0
$data[112,3] = 3
$data[112,4] = 61
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 1
$data[113,0] = 'Bolivia'
$data[113,1] = 61
$data[113,2] = 0
$data[113,3] = 0
$data[113,4] = 61
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 0
$data[114,0] = 'Kirguistan'
$data[114,1] = 58
$data[114,2] = 14
$data[114,3] = 0
$data[114,4] = 58
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 0
$data[115,0] = 'Liechtenstein'
$data[115,1] = 56
$data[115,2] = 0
$data[115,3] = 0
$data[115,4] = 56
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 0
$data[116,0] = 'Gibraltar'
$data[116,1] = 55
$data[116,2] = 20
$data[116,3] = 14
$data[116,4] = 41
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 0
$data[117,0] = 'Paraguay'
$data[117,1] = 52
$data[117,2] = 11
$data[117,3] = 1
$data[117,4] = 48
$data[117,5] = 1
$data[117,6] = 0
$data[117,7] = 3
$data[118,0] = 'Consejo Danes para los Refugiados'
$data[118,1] = 51
$data[118,2] = 0
$data[118,3] = 2
$data[118,4] = 46
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 3
$data[119,0] = 'Mayotte'
$data[119,1] = 50
$data[119,2] = 14
$data[119,3] = 0
$data[119,4] = 50
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = 'Ruanda'
$data[120,1] = 50
$data[120,2] = 0
$data[120,3] = 0
$data[120,4] = 50
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 0
$data[121,0] = 'Banglades'
$data[121,1] = 48
$data[121,2] = 4
$data[121,3] = 11
$data[121,4] = 32
$data[121,5] = 1
$data[121,6] = 0
$data[121,7] = 5
$data[122,0] = 'Puerto Rico'
$data[122,1] = 39
$data[122,2] = 0
$data[122,3] = 1
$data[122,4] = 36
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 2
$data[123,0] = 'Macao'
$data[123,1] = 34
$data[123,2] = 1
$data[123,3] = 10
$data[123,4] = 24
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 0
$data[124,0] = 'Monaco'
$data[124,1] = 33
$data[124,2] = 0
$data[124,3] = 1
$data[124,4] = 32
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = 'Guam'
$data[125,1] = 32
$data[125,2] = 0
$data[125,3] = 0
$data[125,4] = 31
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 1
$data[126,0] = 'Kenia'
$data[126,1] = 31
$data[126,2] = 0
$data[126,3] = 1
$data[126,4] = 29
$data[126,5] = 2
$data[126,6] = 0
$data[126,7] = 1
$data[127,0] = 'Polinesia Francesa'
$data[127,1] = 30
$data[127,2] = 0
$data[127,3] = 0
$data[127,4] = 30
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 0
$data[128,0] = 'Isla de Man'
$data[128,1] = 29
$data[128,2] = 3
$data[128,3] = 0
$data[128,4] = 29
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = 'Aruba'
$data[129,1] = 28
$data[129,2] = 0
$data[129,3] = 1
$data[129,4] = 27
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 0
$data[130,0] = 'Guayana Francesa'
$data[130,1] = 28
$data[130,2] = 0
$data[130,3] = 6
$data[130,4] = 22
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = 'Jamaica'
$data[131,1] = 26
$data[131,2] = 0
$data[131,3] = 2
$data[131,4] = 23
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 1
$data[132,0] = 'Togo'
$data[132,1] = 25
$data[132,2] = 1
$data[132,3] = 1
$data[132,4] = 24
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 0
$data[133,0] = 'Guatemala'
$data[133,1] = 25
$data[133,2] = 0
$data[133,3] = 4
$data[133,4] = 20
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 1
$data[134,0] = 'Barbados'
$data[134,1] = 24
$data[134,2] = 0
$data[134,3] = 0
$data[134,4] = 24
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 0
$data[135,0] = 'Madagascar'
$data[135,1] = 24
$data[135,2] = 1
$data[135,3] = 0
$data[135,4] = 24
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = 'Zambia'
$data[136,1] = 22
$data[136,2] = 6
$data[136,3] = 0
$data[136,4] = 22
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = 'Uganda'
$data[137,1] = 18
$data[137,2] = 4
$data[137,3] = 0
$data[137,4] = 18
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 0
$data[138,0] = 'Islas Virgenes de los Estados Unidos'
$data[138,1] = 17
$data[138,2] = 0
$data[138,3] = 0
$data[138,4] = 17
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = 'Etiopia'
$data[139,1] = 16
$data[139,2] = 4
$data[139,3] = 0
$data[139,4] = 16
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 0
$data[140,0] = 'Nueva Caledonia'
$data[140,1] = 15
$data[140,2] = 1
$data[140,3] = 0
$data[140,4] = 15
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 0
$data[141,0] = 'Bermudas'
$data[141,1] = 15
$data[141,2] = 0
$data[141,3] = 2
$data[141,4] = 13
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = 'Maldivas'
$data[142,1] = 14
$data[142,2] = 1
$data[142,3] = 9
$data[142,4] = 5
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 0
$data[143,0] = 'El Salvador'
$data[143,1] = 13
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 13
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = 'Tanzania'
$data[144,1] = 13
$data[144,2] = 0
$data[144,3] = 1
$data[144,4] = 12
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = 'Guinea Ecuatorial'
$data[145,1] = 12
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 12
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 0
$data[146,0] = 'Republica de Yibuti'
$data[146,1] = 12
$data[146,2] = 1
$data[146,3] = 0
$data[146,4] = 12
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 0
$data[147,0] = 'Mongolia'
$data[147,1] = 11
$data[147,2] = 0
$data[147,3] = 0
$data[147,4] = 11
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 0
$data[148,0] = 'Mali'
$data[148,1] = 11
$data[148,2] = 7
$data[148,3] = 0
$data[148,4] = 11
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 0
$data[149,0] = 'Dominica'
$data[149,1] = 11
$data[149,2] = 0
$data[149,3] = 0
$data[149,4] = 11
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 0
$data[150,0] = 'San Martin (Parte Francesa)'
$data[150,1] = 11
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 11
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 0
$data[151,0] = 'Niger'
$data[151,1] = 10
$data[151,2] = 0
$data[151,3] = 0
$data[151,4] = 9
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 1
$data[152,0] = 'Groenlandia'
$data[152,1] = 10
$data[152,2] = 4
$data[152,3] = 2
$data[152,4] = 8
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 0
$data[153,0] = 'Suazilandia'
$data[153,1] = 9
$data[153,2] = 3
$data[153,3] = 0
$data[153,4] = 9
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = 'Bahamas'
$data[154,1] = 9
$data[154,2] = 0
$data[154,3] = 1
$data[154,4] = 8
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 0
$data[155,0] = 'Surinam'
$data[155,1] = 8
$data[155,2] = 0
$data[155,3] = 0
$data[155,4] = 8
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = 'Haiti'
$data[156,1] = 8
$data[156,2] = 0
$data[156,3] = 0
$data[156,4] = 8
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = 'Guinea'
$data[157,1] = 8
$data[157,2] = 4
$data[157,3] = 0
$data[157,4] = 8
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 0
$data[158,0] = 'Islas Caimanes'
$data[158,1] = 8
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 7
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 1
$data[159,0] = 'Namibia'
$data[159,1] = 8
$data[159,2] = 0
$data[159,3] = 2
$data[159,4] = 6
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = 'Mozambique'
$data[160,1] = 7
$data[160,2] = 0
$data[160,3] = 0
$data[160,4] = 7
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = 'Antigua y Barbuda'
$data[161,1] = 7
$data[161,2] = 0
$data[161,3] = 0
$data[161,4] = 7
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = 'Granada'
$data[162,1] = 7
$data[162,2] = 0
$data[162,3] = 0
$data[162,4] = 7
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = 'Seychelles'
$data[163,1] = 7
$data[163,2] = 0
$data[163,3] = 0
$data[163,4] = 7
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = 'Gabon'
$data[164,1] = 7
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 6
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 1
$data[165,0] = 'Curazao'
$data[165,1] = 7
$data[165,2] = 0
$data[165,3] = 2
$data[165,4] = 4
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 1
$data[166,0] = 'Eritrea'
$data[166,1] = 6
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 6
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = 'Laos'
$data[167,1] = 6
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 6
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = 'Benin'
$data[168,1] = 6
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 6
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = 'San Bartolome'
$data[169,1] = 5
$data[169,2] = 2
$data[169,3] = 0
$data[169,4] = 5
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = 'Siria'
$data[170,1] = 5
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 5
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = 'Birmania'
$data[171,1] = 5
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 5
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = 'Montserrat'
$data[172,1] = 5
$data[172,2] = 0
$data[172,3] = 0
$data[172,4] = 5
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = 'Fiyi'
$data[173,1] = 5
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 5
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = 'Cabo Verde'
$data[174,1] = 5
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 4
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 1
$data[175,0] = 'Guyana'
$data[175,1] = 5
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 4
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 1
$data[176,0] = 'Zimbabue'
$data[176,1] = 5
$data[176,2] = 2
$data[176,3] = 0
$data[176,4] = 4
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 1
$data[177,0] = 'Angola'
$data[177,1] = 4
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 4
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Congo'
$data[178,1] = 4
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 4
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = 'Santa Sede'
$data[179,1] = 4
$data[179,2] = 0
$data[179,3] = 0
$data[179,4] = 4
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = 'Nepal'
$data[180,1] = 4
$data[180,2] = 1
$data[180,3] = 1
$data[180,4] = 3
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = 'Republica de Africa Central'
$data[181,1] = 3
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 3
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = 'San Martin (Parte Holandesa)'
$data[182,1] = 3
$data[182,2] = 0
$data[182,3] = 0
$data[182,4] = 3
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Somalia'
$data[183,1] = 3
$data[183,2] = 1
$data[183,3] = 0
$data[183,4] = 3
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Butan'
$data[184,1] = 3
$data[184,2] = 1
$data[184,3] = 0
$data[184,4] = 3
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Mauritania'
$data[185,1] = 3
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 3
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = 'Liberia'
$data[186,1] = 3
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 3
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 0
$data[187,0] = 'Republica del Chad'
$data[187,1] = 3
$data[187,2] = 0
$data[187,3] = 0
$data[187,4] = 3
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = 'Gambia'
$data[188,1] = 3
$data[188,2] = 0
$data[188,3] = 0
$data[188,4] = 2
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 1
$data[189,0] = 'Sudan'
$data[189,1] = 3
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 2
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 1
$data[190,0] = 'Santa Lucia'
$data[190,1] = 3
$data[190,2] = 0
$data[190,3] = 1
$data[190,4] = 2
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = 'Anguila'
$data[191,1] = 2
$data[191,2] = 0
$data[191,3] = 0
$data[191,4] = 2
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = 'Guinea-Bisau'
$data[192,1] = 2
$data[192,2] = 0
$data[192,3] = 0
$data[192,4] = 2
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = 'San Cristobal y Nieves'
$data[193,1] = 2
$data[193,2] = 0
$data[193,3] = 0
$data[193,4] = 2
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = 'Islas Virgenes Britanicas'
$data[194,1] = 2
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 2
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 0
$data[195,0] = 'Belice'
$data[195,1] = 2
$data[195,2] = 0
$data[195,3] = 0
$data[195,4] = 2
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$data[196,0] = 'Islas Turcas y Caicos'
$data[196,1] = 2
$data[196,2] = 0
$data[196,3] = 0
$data[196,4] = 2
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = 'Nicaragua'
$data[197,1] = 2
$data[197,2] = 0
$data[197,3] = 0
$data[197,4] = 1
$data[197,5] = 0
$data[197,6] = 1
$data[197,7] = 1
$data[198,0] = 'San Vicente y las Granadinas'
$data[198,1] = 1
$data[198,2] = 0
$data[198,3] = 0
$data[198,4] = 1
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = 'Papua Nueva Guinea'
$data[199,1] = 1
$data[199,2] = 0
$data[199,3] = 0
$data[199,4] = 1
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = 'Libia'
$data[200,1] = 1
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 1
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = 'Timor Oriental'
$data[201,1] = 1
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 1
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0

$ws.Range("A4:H205").Value2 = $data

